# Updated cryptos list on Tue Nov 28 10:15:13 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.056.77"
$ws.Range("E2").Value = "  -0.80%  "

$ws.Range("D3").Value = "2.017.03"
$ws.Range("E3").Value = "  -1.36%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.89"
$ws.Range("E5").Value = "  -0.40%  "

$ws.Range("E6").Value = "  -0.91%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.56"
$ws.Range("E8").Value = "  -1.12%  "

$ws.Range("E9").Value = "  -2.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0778"
$ws.Range("E10").Value = "  -3.21%  "

$ws.Range("E11").Value = "  -4.08%  "

$ws.Range("D12").Value = "2.314.53"
$ws.Range("E12").Value = "  -1.42%  "

$ws.Range("E13").Value = "  -2.17%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.91"
$ws.Range("E14").Value = "  -3.47%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.19"
$ws.Range("E15").Value = "  -1.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.737"
$ws.Range("E16").Value = "  -2.15%  "

$ws.Range("D17").Value = "2.017.06"
$ws.Range("E17").Value = "  -1.53%  "

$ws.Range("D18").Value = "37.017.38"
$ws.Range("E18").Value = "  -0.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.10"
$ws.Range("E19").Value = "  +0.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.89"
$ws.Range("E20").Value = "  -0.95%  "

$ws.Range("D21").Value = "0.0₃0813"
$ws.Range("E21").Value = "  -4.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.27"
$ws.Range("E22").Value = "  -0.88%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.13%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").Value = "  +2.30%  "

$ws.Range("E25").Value = "  -4.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.44"
$ws.Range("E26").Value = "  -2.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.99"
$ws.Range("E27").Value = "  -5.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.126"
$ws.Range("E28").Value = "  -0.50%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.64"
$ws.Range("E29").Value = "  -1.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.30"
$ws.Range("E30").Value = "  -3.76%  "

$ws.Range("E31").Value = "  -0.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.42"
$ws.Range("E32").Value = "  -2.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0600"
$ws.Range("E33").Value = "  -2.00%  "

$ws.Range("E34").Value = "  -1.42%  "

$ws.Range("E35").Value = "  -0.84%  "

$ws.Range("E36").Value = "  +2.03%  "

$ws.Range("E37").Value = "  +0.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.15"
$ws.Range("E38").Value = "  -1.28%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.47"
$ws.Range("E39").Value = "  +0.75%  "

$ws.Range("D40").Value = "1.469.30"
$ws.Range("E40").Value = "  -2.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0212"
$ws.Range("E41").Value = "  -3.58%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "94.25"
$ws.Range("E42").Value = "  -1.58%  "

$ws.Range("E43").Value = "  -3.35%  "

$ws.Range("E44").Value = "  -3.66%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.07"
$ws.Range("E45").Value = "  -3.30%  "

$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.11"
$ws.Range("E46").Value = "  -2.83%  "

$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.10"
$ws.Range("E47").Value = "  +8.54%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.998"
$ws.Range("E48").Value = "  -1.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.05"
$ws.Range("E49").Value = "  -1.97%  "

$ws.Range("E50").Value = "  -0.06%  "

$ws.Range("D51").Value = "2.204.38"
$ws.Range("E51").Value = "  -1.44%  "
